# Update Metadata last-updated timestamp
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:43 AM"

# Insert a new top data row in Stock List, shifting existing rows down,
# and drop the row that falls off the bottom of the table (row 77).
$ws = $wb.Worksheets.Item("Stock List")
$ws.Rows.Item(2).Insert()
$ws.Range("A2:H2").ClearFormats()

$ws.Range("A2").Value = "📋"
$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 0

$ws.Rows.Item(77).Delete()
